# Massive updates to monsoon variability
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the monthly adjustment rows (Je, Jc, Jp, Jn -> rows 14-17)
$ws.Range("B14").Value = 0.1
$ws.Range("D14").Value = 0.75

$ws.Range("B15").Value = 0.1
$ws.Range("D15").Value = 0.75

$ws.Range("B16").Value = 0.1
$ws.Range("D16").Value = 0.75

$ws.Range("B17").Value = 0.1
$ws.Range("D17").Value = 0.75

# Match the active selection captured in the saved file
$ws.Range("D15").Select()
